$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain value / string / number updates ---
$ws.Range("A2").Value = "Martin Guptill"
$ws.Range("B2").Value = 45
$ws.Range("C2").Value = 13
$ws.Range("D2").Value = "LBW"
$ws.Range("E2").Value = " Hardik Pandya"
$ws.Range("J2").Value = "KL Rahul"
$ws.Range("K2").Value = 6
$ws.Range("L2").Value = 4
$ws.Range("M2").Value = "LBW"
$ws.Range("N2").Value = " Trent Boult"
$ws.Range("A3").Value = "Daryl Mitchell"
$ws.Range("B3").Value = 33
$ws.Range("C3").Value = 11
$ws.Range("E3").Value = " Bhuvneshwar Kumar"
$ws.Range("J3").Value = "Rohit Sharma"
$ws.Range("K3").Value = 12
$ws.Range("L3").Value = 4
$ws.Range("N3").Value = " Tim Southee"
$ws.Range("A4").Value = "Kane Williamson(C)"
$ws.Range("B4").Value = 11
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = "Caught"
$ws.Range("E4").Value = " Hardik Pandya"
$ws.Range("J4").Value = "Virat Kohli(C)"
$ws.Range("K4").Value = 20
$ws.Range("L4").Value = 7
$ws.Range("M4").Value = "LBW"
$ws.Range("N4").Value = " Adam Milne"
$ws.Range("A5").Value = "Devon Conway"
$ws.Range("B5").Value = 43
$ws.Range("C5").Value = 13
$ws.Range("E5").Value = " Mohommad Shami"
$ws.Range("J5").Value = "Suryakumar Yadav"
$ws.Range("K5").Value = 11
$ws.Range("L5").Value = 5
$ws.Range("M5").Value = "Caught"
$ws.Range("N5").Value = " Ish Sodhi"
$ws.Range("A6").Value = "Glenn Phillips"
$ws.Range("B6").Value = 23
$ws.Range("C6").Value = 6
$ws.Range("E6").Value = " Jasprit Bumrah"
$ws.Range("J6").Value = "Rishabh Pant"
$ws.Range("K6").Value = 6
$ws.Range("L6").Value = 4
$ws.Range("M6").Value = "NOT OUT"
$ws.Range("N6").Value = " "
$ws.Range("A7").Value = "James Neesham"
$ws.Range("B7").Value = 14
$ws.Range("C7").Value = 4
$ws.Range("E7").Value = " Hardik Pandya"
$ws.Range("J7").Value = "Ravindra Jadeja"
$ws.Range("K7").Value = 12
$ws.Range("L7").Value = 3
$ws.Range("M7").Value = "Bowled"
$ws.Range("N7").Value = " Mitchell Santner"
$ws.Range("A8").Value = "Mitchell Santner"
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = "Caught"
$ws.Range("E8").Value = " Bhuvneshwar Kumar"
$ws.Range("J8").Value = "Hardik Pandya"
$ws.Range("K8").Value = 12
$ws.Range("L8").Value = 4
$ws.Range("M8").Value = "LBW"
$ws.Range("N8").Value = " Mitchell Santner"
$ws.Range("A9").Value = "Adam Milne"
$ws.Range("B9").Value = 8
$ws.Range("C9").Value = 6
$ws.Range("D9").Value = "LBW"
$ws.Range("E9").Value = " Jasprit Bumrah"
$ws.Range("J9").Value = "Bhuvneshwar Kumar"
$ws.Range("K9").Value = 8
$ws.Range("L9").Value = 5
$ws.Range("M9").Value = "Bowled"
$ws.Range("N9").Value = " Trent Boult"
$ws.Range("A10").Value = "Ish Sodhi"
$ws.Range("B10").Value = 20
$ws.Range("C10").Value = 10
$ws.Range("D10").Value = "Bowled"
$ws.Range("E10").Value = " Kuldeep Yadav"
$ws.Range("J10").Value = "Mohommad Shami"
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 1
$ws.Range("N10").Value = " Trent Boult"
$ws.Range("A11").Value = "Tim Southee"
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = "NOT OUT"
$ws.Range("E11").Value = " "
$ws.Range("J11").Value = "Jasprit Bumrah"
$ws.Range("K11").Value = 5
$ws.Range("M11").Value = "LBW"
$ws.Range("N11").Value = " Tim Southee"
$ws.Range("A12").Value = "Trent Boult"
$ws.Range("B12").Value = 5
$ws.Range("C12").Value = 3
$ws.Range("D12").Value = "Bowled"
$ws.Range("E12").Value = " Jasprit Bumrah"
$ws.Range("J12").Value = "Kuldeep Yadav"
$ws.Range("K12").Value = 1
$ws.Range("N12").Value = " Tim Southee"
$ws.Range("A16").Value = 205
$ws.Range("D16").Value = 72
$ws.Range("J16").Value = 93
$ws.Range("M16").Value = 41
$ws.Range("A21").Value = "Mohommad Shami"
$ws.Range("C21").Value = 38
$ws.Range("D21").Value = 1
$ws.Range("J21").Value = "Ish Sodhi"
$ws.Range("L21").Value = 14
$ws.Range("N21").Value = 14
$ws.Range("A22").Value = "Bhuvneshwar Kumar"
$ws.Range("C22").Value = 37
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = 18.5
$ws.Range("J22").Value = "Adam Milne"
$ws.Range("L22").Value = 23
$ws.Range("M22").Value = 1
$ws.Range("N22").Value = 23
$ws.Range("A23").Value = "Hardik Pandya"
$ws.Range("C23").Value = 19
$ws.Range("E23").Value = 9.5
$ws.Range("J23").Value = "Mitchell Santner"
$ws.Range("L23").Value = 12
$ws.Range("M23").Value = 2
$ws.Range("N23").Value = 12
$ws.Range("A24").Value = "Kuldeep Yadav"
$ws.Range("C24").Value = 62
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 20.67
$ws.Range("J24").Value = "Trent Boult"
$ws.Range("A25").Value = "Jasprit Bumrah"
$ws.Range("C25").Value = 49
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 16.33
$ws.Range("J25").Value = "Tim Southee"
$ws.Range("M25").Value = 3
$ws.Range("N25").Value = 12

# --- Numeric-looking text cells (overs, e.g. "2.0") must stay text, not become numbers ---
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "12.0"
$ws.Range("C16").Style = "Normal"
$ws.Range("L16").NumberFormat = "@"
$ws.Range("L16").Value = "6.5"
$ws.Range("L16").Style = "Normal"
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "2.0"
$ws.Range("B21").Style = "Normal"
$ws.Range("K21").NumberFormat = "@"
$ws.Range("K21").Value = "2.0"
$ws.Range("K21").Style = "Normal"
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "2.0"
$ws.Range("B22").Style = "Normal"
$ws.Range("K22").NumberFormat = "@"
$ws.Range("K22").Value = "1.0"
$ws.Range("K22").Style = "Normal"
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "2.0"
$ws.Range("B23").Style = "Normal"
$ws.Range("K23").NumberFormat = "@"
$ws.Range("K23").Value = "1.0"
$ws.Range("K23").Style = "Normal"
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "2.0"
$ws.Range("B24").Style = "Normal"
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "3.0"
$ws.Range("B25").Style = "Normal"
$ws.Range("K25").NumberFormat = "@"
$ws.Range("K25").Value = "1.5"
$ws.Range("K25").Style = "Normal"
